$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("report")

$ws.Range("F4").Value = "Was not able to save the task of CPI."
$ws.Range("G4").Value = "2022-09-08 19:05:04"
